$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39578.5
$ws.Range("J3").Value = 39578.5
$ws.Range("L3").Value = 39578.5
$ws.Range("N3").Value = -39806.5
$ws.Range("H7").Value = 34899.5
$ws.Range("J7").Value = 34899.5
$ws.Range("L7").Value = 34899.5
$ws.Range("N7").Value = -35123.5
$ws.Range("H14").Value = 34899.5
$ws.Range("J14").Value = 34899.5
$ws.Range("L14").Value = 34899.5
$ws.Range("N14").Value = -35281.5
$ws.Range("H19").Value = 2925630
$ws.Range("I19").Value = 5264297.5
$ws.Range("K19").Value = 5264297.5
$ws.Range("M19").Value = -5264122.5
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H96").Value = 1371
$ws.Range("I96").Value = 1371
$ws.Range("K96").Value = 4113
$ws.Range("M96").Value = -2740
$ws.Range("H102").Value = 39578.5
$ws.Range("J102").Value = 39578.5
$ws.Range("L102").Value = 39578.5
$ws.Range("N102").Value = -46068.5
$ws.Range("H107").Value = 1907.8125
$ws.Range("I107").Value = 1901.9231
$ws.Range("J107").Value = 1933.3334
$ws.Range("K107").Value = 1901.9231
$ws.Range("L107").Value = 1933.3334
$ws.Range("M107").Value = 18.07690000000002
$ws.Range("N107").Value = -5773.3334
$ws.Range("H113").Value = 13092.818
$ws.Range("I113").Value = 2314.6667
$ws.Range("J113").Value = 17134.625
$ws.Range("K113").Value = 2314.6667
$ws.Range("L113").Value = 17134.625
$ws.Range("M113").Value = 939.3332999999998
$ws.Range("N113").Value = -23642.625
$ws.Range("H116").Value = 5513.857
$ws.Range("I116").Value = 2095.3333
$ws.Range("J116").Value = 7133.1577
$ws.Range("K116").Value = 2095.3333
$ws.Range("L116").Value = 7133.1577
$ws.Range("M116").Value = 1346.6667
$ws.Range("N116").Value = -14017.1577

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4234.8335
$ws.Range("I25").Value = 1607.25
$ws.Range("J25").Value = 9490
$ws.Range("K25").Value = 1607.25
$ws.Range("L25").Value = 9490
$ws.Range("M25").Value = -1205.25
$ws.Range("N25").Value = -10294
$ws.Range("H32").Value = 5677.392
$ws.Range("I32").Value = 4349.1704
$ws.Range("K32").Value = 4349.1704
$ws.Range("M32").Value = -4062.1704
$ws.Range("H88").Value = 6669153
$ws.Range("I88").Value = 13334929
$ws.Range("J88").Value = 3377.4
$ws.Range("K88").Value = 13334929
$ws.Range("L88").Value = 3377.4
$ws.Range("M88").Value = -13334523
$ws.Range("N88").Value = -4189.4
$ws.Range("H91").Value = 6669153
$ws.Range("I91").Value = 13334929
$ws.Range("J91").Value = 3377.4
$ws.Range("K91").Value = 13334929
$ws.Range("L91").Value = 3377.4
$ws.Range("M91").Value = -13333525
$ws.Range("N91").Value = -6185.4
$ws.Range("H102").Value = 2002.75
$ws.Range("I102").Value = 2170
$ws.Range("J102").Value = 1501
$ws.Range("K102").Value = 2170
$ws.Range("L102").Value = 1501
$ws.Range("M102").Value = -548
$ws.Range("N102").Value = -4745
$ws.Range("H132").Value = 2786.2778
$ws.Range("I132").Value = 1215.1
$ws.Range("K132").Value = 3645.3
$ws.Range("M132").Value = -1115.3

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 32920.5
$ws.Range("J29").Value = 32920.5
$ws.Range("L29").Value = 32920.5
$ws.Range("N29").Value = -33506.5
$ws.Range("H62").Value = 3266.8333
$ws.Range("I62").Value = 3001.6667
$ws.Range("J62").Value = 3532
$ws.Range("K62").Value = 3001.6667
$ws.Range("L62").Value = 3532
$ws.Range("M62").Value = -2377.6667
$ws.Range("N62").Value = -4780
$ws.Range("H65").Value = 3266.8333
$ws.Range("I65").Value = 3001.6667
$ws.Range("J65").Value = 3532
$ws.Range("K65").Value = 15008.3335
$ws.Range("L65").Value = 17660
$ws.Range("M65").Value = -11888.3335
$ws.Range("N65").Value = -23900

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 83090.91
$ws.Range("J37").Value = 83090.91
$ws.Range("L37").Value = 249272.73
$ws.Range("N37").Value = -249496.73

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 29490.25
$ws.Range("J96").Value = 29490.25
$ws.Range("L96").Value = 29490.25
$ws.Range("N96").Value = -34982.25

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6422.222
$ws.Range("I132").Value = 2211.8333
$ws.Range("J132").Value = 9790.532999999999
$ws.Range("K132").Value = 6635.499899999999
$ws.Range("L132").Value = 29371.599
$ws.Range("M132").Value = -4105.499899999999
$ws.Range("N132").Value = -34431.599
$ws.Range("H136").Value = 4673.136
$ws.Range("I136").Value = 1529.8572
$ws.Range("K136").Value = 4589.571599999999
$ws.Range("M136").Value = -2039.571599999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
$ws.Range("H34").Value = 41014.5
$ws.Range("I34").Value = 12000
$ws.Range("J34").Value = 70029
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 70029
$ws.Range("M34").Value = -11797
$ws.Range("N34").Value = -70435
$ws.Range("H109").Value = 28377
$ws.Range("J109").Value = 28377
$ws.Range("L109").Value = 28377
$ws.Range("N109").Value = -31151
$ws.Range("H113").Value = 4551.7827
$ws.Range("I113").Value = 7899.846
$ws.Range("J113").Value = 199.3
$ws.Range("K113").Value = 23699.538
$ws.Range("L113").Value = 597.9000000000001
$ws.Range("M113").Value = -21529.538
$ws.Range("N113").Value = -4937.9
$ws.Range("H132").Value = 27789278
$ws.Range("I132").Value = 26199.75
$ws.Range("K132").Value = 78599.25
$ws.Range("M132").Value = -76069.25
